$wb = $excel.ActiveWorkbook

# --- Rename the first sheet (Ex1 -> ex1); content/data unchanged ---
$ws1 = $wb.Worksheets.Item("Ex1")
$ws1.Name = "ex1"

# --- Build the replacement "ex2" sheet so it gets a fresh sheetId (3) ---
# Duplicate the "ex1" tab (placed right after "Ex2") while the old "Ex2"
# (sheetId 2) still exists, so the duplicate is allocated sheetId 3; then
# delete the old "Ex2" tab, rename the duplicate into its place, and wipe
# its (copied) contents so it can be repopulated from scratch below.
$oldEx2 = $wb.Worksheets.Item("Ex2")
$ws1.Copy($null, $oldEx2)
[void]$oldEx2.Delete()

# Re-fetch: the handle captured before the Delete() above goes stale once the
# sibling sheet is removed, so further writes through it are silently dropped.
$newEx2 = $wb.Worksheets.Item("ex1 (2)")
$newEx2.Name = "ex2"
$newEx2.Cells.Clear()

# --- Populate the new "ex2" sheet with the Contractual Setting example data ---
# Write header cells in this order so the shared-strings table comes out in
# the same sequence as the target workbook (activeCust, then Period, then lostCust).
$newEx2.Range("B1").Value = "activeCust"
$newEx2.Range("A1").Value = "Period"
$newEx2.Range("C1").Value = "lostCust"

$data = @(
    @(1, 869, 131),
    @(2, 743, 126),
    @(3, 653, 90),
    @(4, 593, 60),
    @(5, 551, 42),
    @(6, 517, 34),
    @(7, 491, 26)
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $newEx2.Cells.Item($row, 1).Value = $data[$i][0]
    $newEx2.Cells.Item($row, 2).Value = $data[$i][1]
    $newEx2.Cells.Item($row, 3).Value = $data[$i][2]
}

# --- Selection / active-tab state matches the target file ---
[void]$ws1.Range("B8").Select()

$newEx2.Activate()
[void]$newEx2.Range("A9:XFD17").Select()
